# Fixing geopoint in shared_table model
# Insert a new "schema.name" column into the "model" sheet (right after
# "schema.type", before "schema.elementType"), populate it for the
# refrigerator_location (geopoint) row, and rename the
# schema.properties.* headers to schema.properties.*.type.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("model")

# Insert a new column before column C (schema.elementType), shifting
# schema.elementType and the schema.properties.* columns one to the right.
$ws.Columns("C:C").Insert()

# New header for the inserted column.
$ws.Cells.Item(1, 3).Value = "schema.name"

# Rename the schema.properties.* headers (now shifted to columns E-H) to
# include the ".type" suffix.
$ws.Cells.Item(1, 5).Value = "schema.properties.latitude.type"
$ws.Cells.Item(1, 6).Value = "schema.properties.longitude.type"
$ws.Cells.Item(1, 7).Value = "schema.properties.altitude.type"
$ws.Cells.Item(1, 8).Value = "schema.properties.accuracy.type"

# Populate the new schema.name value for the refrigerator_location
# (geopoint) row - same as its elementType value.
$ws.Cells.Item(4, 3).Value = "geopoint"
